# Fix ingest perf decrease after switching DB engines
# Adds a new "postgres switch" results column (H) to the perf-stats sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header ---
$ws.Range("H1").Value = "postgres switch"

# --- New column data (values copied from the new DB-engine benchmark run) ---
$ws.Range("H2").Value  = 5000000
$ws.Range("H3").Value  = 1334.82
$ws.Range("H4").Value  = 2615660
$ws.Range("H5").Value  = 341491
$ws.Range("H6").Value  = 69696
$ws.Range("H7").Value  = 8.26
$ws.Range("H8").Value  = 1921455
$ws.Range("H9").Value  = 336.06
$ws.Range("H10").Value = 450.99
$ws.Range("H11").Value = 25.06

# --- Column widths: widen G slightly, size the new H column ---
$ws.Columns.Item(7).ColumnWidth = 14.6
$ws.Columns.Item(8).ColumnWidth = 12.6

# --- Selection moves to the cell below the new column's last entry ---
$ws.Range("H12").Select() | Out-Null
